$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("D2").Value = 44417
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 1300
$ws.Range("O2").Value = 1400
$ws.Range("P2").Value = 1350
$ws.Range("S2").Value = 1350

# Row 3
$ws.Range("D3").Value = 44372
$ws.Range("M3").Value = 900
$ws.Range("N3").Value = 750
$ws.Range("O3").Value = 800
$ws.Range("P3").Value = 772
$ws.Range("S3").Value = 772

# Row 4
$ws.Range("D4").Value = 44372
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 900
$ws.Range("N4").Value = 600
$ws.Range("O4").Value = 650
$ws.Range("P4").Value = 628
$ws.Range("S4").Value = 628

# Row 5
$ws.Range("D5").Value = 44260
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 1900
$ws.Range("O5").Value = 2000
$ws.Range("P5").Value = 1950
$ws.Range("S5").Value = 1950

# Row 6
$ws.Range("D6").Value = 44351
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 700
$ws.Range("O6").Value = 800
$ws.Range("P6").Value = 750
$ws.Range("S6").Value = 750

# Row 7
$ws.Range("D7").Value = 44351
$ws.Range("N7").Value = 600
$ws.Range("O7").Value = 700
$ws.Range("P7").Value = 650
$ws.Range("S7").Value = 650

# Row 8
$ws.Range("D8").Value = 44358
$ws.Range("M8").Value = 200

# Row 9
$ws.Range("D9").Value = 44358
$ws.Range("M9").Value = 200
$ws.Range("O9").Value = 650
$ws.Range("P9").Value = 625
$ws.Range("S9").Value = 625

# Row 10
$ws.Range("D10").Value = 44326
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 600
$ws.Range("O10").Value = 700
$ws.Range("P10").Value = 650
$ws.Range("S10").Value = 650

# Row 11
$ws.Range("D11").Value = 44407
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 200
$ws.Range("O11").Value = 650
$ws.Range("P11").Value = 625
$ws.Range("S11").Value = 625

# Row 14
$ws.Range("D14").Value = 44316
$ws.Range("M14").Value = 140
$ws.Range("N14").Value = 1100
$ws.Range("O14").Value = 1200
$ws.Range("P14").Value = 1150
$ws.Range("S14").Value = 1150

# Row 15
$ws.Range("D15").Value = 44330
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 1200
$ws.Range("O15").Value = 1300
$ws.Range("P15").Value = 1250
$ws.Range("S15").Value = 1250

# Row 16
$ws.Range("D16").Value = 44330
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 1000
$ws.Range("O16").Value = 1100
$ws.Range("P16").Value = 1050
$ws.Range("S16").Value = 1050

# Row 17
$ws.Range("D17").Value = 44403
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 1200
$ws.Range("O17").Value = 1300
$ws.Range("P17").Value = 1250
$ws.Range("S17").Value = 1250

# Row 18
$ws.Range("D18").Value = 44403
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 120
$ws.Range("N18").Value = 950
$ws.Range("O18").Value = 1000
$ws.Range("P18").Value = 975
$ws.Range("S18").Value = 975

# Row 19
$ws.Range("D19").Value = 44379
$ws.Range("M19").Value = 150
$ws.Range("N19").Value = 700
$ws.Range("O19").Value = 800
$ws.Range("P19").Value = 747
$ws.Range("S19").Value = 747

# Row 20
$ws.Range("D20").Value = 44379
$ws.Range("M20").Value = 140
$ws.Range("N20").Value = 500
$ws.Range("O20").Value = 600
$ws.Range("P20").Value = 543
$ws.Range("S20").Value = 543

# Row 21
$ws.Range("D21").Value = 44498
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 1200
$ws.Range("O21").Value = 1300
$ws.Range("P21").Value = 1250
$ws.Range("S21").Value = 1250

# Row 22
$ws.Range("D22").Value = 44425
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 140
$ws.Range("N22").Value = 1200
$ws.Range("O22").Value = 1300
$ws.Range("P22").Value = 1250
$ws.Range("S22").Value = 1250

# Row 23
$ws.Range("D23").Value = 44389
$ws.Range("M23").Value = 140
$ws.Range("N23").Value = 750
$ws.Range("O23").Value = 800
$ws.Range("P23").Value = 775
$ws.Range("S23").Value = 775

# Row 24
$ws.Range("D24").Value = 44389
$ws.Range("L24").Value = "Segunda"
$ws.Range("M24").Value = 120
$ws.Range("N24").Value = 600
$ws.Range("O24").Value = 700
$ws.Range("P24").Value = 650
$ws.Range("S24").Value = 650

# Row 25
$ws.Range("D25").Value = 44386
$ws.Range("M25").Value = 160
$ws.Range("N25").Value = 700
$ws.Range("O25").Value = 750
$ws.Range("P25").Value = 725
$ws.Range("S25").Value = 725

# Row 26
$ws.Range("D26").Value = 44386
$ws.Range("M26").Value = 200
$ws.Range("P26").Value = 625
$ws.Range("S26").Value = 625

# Row 27
$ws.Range("D27").Value = 44414
$ws.Range("M27").Value = 160
$ws.Range("N27").Value = 1300
$ws.Range("O27").Value = 1400
$ws.Range("P27").Value = 1350
$ws.Range("S27").Value = 1350

# Row 28
$ws.Range("D28").Value = 44350
$ws.Range("M28").Value = 140
$ws.Range("N28").Value = 750
$ws.Range("P28").Value = 775
$ws.Range("S28").Value = 775

# Row 29
$ws.Range("D29").Value = 44473
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 160
$ws.Range("N29").Value = 1500
$ws.Range("O29").Value = 1600
$ws.Range("P29").Value = 1550
$ws.Range("S29").Value = 1550
